$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out rows 5-17 (their card text is being folded into the consolidated
# tuple-strings now written into rows 2-4), shrinking the used range.
$ws.Range("A5:A17").ClearContents() | Out-Null

# Rewrite rows 2-4 with the new consolidated "(name, [fields...])" strings.
$ws.Range("A2").Value = "('Fiery Temper', ['{1}{R}{R}', 'Instant', 'Fiery Temper deals 3 damage to any target.', 'Madness {R} (If you discard this card, discard it into exile. When you do, cast it for its madness cost or put it into your graveyard.)'])"
$ws.Range("A3").Value = "('Icatian Javelineers', ['{W}', 'Creature — Human Soldier', 'Icatian Javelineers enters the battlefield with a javelin counter on it.', '{T}, Remove a javelin counter from Icatian Javelineers: It deals 1 damage to any target.', '1/1'])"
$ws.Range("A4").Value = "('Wood Elves', ['{2}{G}', 'Creature — Elf Scout', 'When Wood Elves enters the battlefield, search your library for a Forest card and put that card onto the battlefield. Then shuffle your library.', '1/1'])"
